$wb = $excel.ActiveWorkbook

# 1. Update status text "Ready for handoff" -> "Handed back: in sync with en-US"
#    This shared string is used on the Overview sheet (columns B "zh-cn" and
#    C "de-de") as well as on the per-language sheets (column C "Status").
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = "Handed back: in sync with en-US"
$ov.Range("C2").Value = "Handed back: in sync with en-US"
$ov.Range("B3").Value = "Handed back: in sync with en-US"
$ov.Range("C3").Value = "Handed back: in sync with en-US"

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("C2").Value = "Handed back: in sync with en-US"
    $ws.Range("C3").Value = "Handed back: in sync with en-US"
}

# 2. zh-cn sheet: add Latest Target File (F) / Latest Handback File (G) columns
#    for the two rows (as real hyperlinks, same as the existing columns
#    A/B/D), and set the handback datetime (H) now that the file is handed
#    back (row 3 shares the same handback-datetime string as row 2).
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/d19d8ad761d2f76a111bc738f2f1ac5a44a672a1/e2e/11be02ce-bfed-43d8-bfaf-d42d007378a1.md", "", "", "11be02ce-bfed-43d8-bfaf-d42d007378a1.md")
$zh.Hyperlinks.Add($zh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/38ce08747f8a4f5be03af76f7c892b46a565fa54/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/11be02ce-bfed-43d8-bfaf-d42d007378a1.eafb28aeb9a246dfda08131c2a4ba31ab61a643b.zh-cn.xlf", "", "", "11be02ce-bfed-43d8-bfaf-d42d007378a1.eafb28aeb9a246dfda08131c2a4ba31ab61a643b.zh-cn.xlf")
$zh.Range("H2").Value = "2016-03-18 05:10:28"
$zh.Range("H3").Value = "2016-03-18 05:10:28"

$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/d19d8ad761d2f76a111bc738f2f1ac5a44a672a1/e2e/83019e79-b3a6-4a13-b49e-3272deba477e.md", "", "", "83019e79-b3a6-4a13-b49e-3272deba477e.md")
$zh.Hyperlinks.Add($zh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/38ce08747f8a4f5be03af76f7c892b46a565fa54/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/83019e79-b3a6-4a13-b49e-3272deba477e.0dc12a42b579df1f7aeab138f87f5d5fec573e6c.zh-cn.xlf", "", "", "83019e79-b3a6-4a13-b49e-3272deba477e.0dc12a42b579df1f7aeab138f87f5d5fec573e6c.zh-cn.xlf")

# 3. de-de sheet: same treatment
$de = $wb.Worksheets.Item("de-de")

$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/d19d8ad761d2f76a111bc738f2f1ac5a44a672a1/e2e/11be02ce-bfed-43d8-bfaf-d42d007378a1.md", "", "", "11be02ce-bfed-43d8-bfaf-d42d007378a1.md")
$de.Hyperlinks.Add($de.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ef8f709527fd789561fd309199a6d2ae1e5abbe7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/11be02ce-bfed-43d8-bfaf-d42d007378a1.eafb28aeb9a246dfda08131c2a4ba31ab61a643b.de-de.xlf", "", "", "11be02ce-bfed-43d8-bfaf-d42d007378a1.eafb28aeb9a246dfda08131c2a4ba31ab61a643b.de-de.xlf")
$de.Range("H2").Value = "2016-03-18 05:10:33"
$de.Range("H3").Value = "2016-03-18 05:10:33"

$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/d19d8ad761d2f76a111bc738f2f1ac5a44a672a1/e2e/83019e79-b3a6-4a13-b49e-3272deba477e.md", "", "", "83019e79-b3a6-4a13-b49e-3272deba477e.md")
$de.Hyperlinks.Add($de.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ef8f709527fd789561fd309199a6d2ae1e5abbe7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/83019e79-b3a6-4a13-b49e-3272deba477e.0dc12a42b579df1f7aeab138f87f5d5fec573e6c.de-de.xlf", "", "", "83019e79-b3a6-4a13-b49e-3272deba477e.0dc12a42b579df1f7aeab138f87f5d5fec573e6c.de-de.xlf")
